# Edit slide 8 ("Collaborations") content placeholder to match the target revision.
$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(8)
$sh = $s.Shapes.Item(2)          # "Content Placeholder 2"
$tr = $sh.TextFrame.TextRange

function Set-ParagraphText {
    param($TextRange, [int]$Index, [string]$NewText)
    $para = $TextRange.Paragraphs($Index, 1)
    # Re-select the paragraph's own characters (not the Paragraphs() wrapper) so the
    # whole run is replaced atomically instead of being diffed/split into two runs.
    $whole = $TextRange.Characters($para.Start, $para.Length)
    $whole.Text = $NewText
}

# 1) "Correlate graph metrics (node degree, node/edge properties) with time-series trends"
#    -> "Correlate graph metrics with time-series trends"
Set-ParagraphText $tr 13 "Correlate graph metrics with time-series trends"

# 2) "...correlation between pressure measurements and dynamic edge weights between nearby sensors"
#    -> "...correlation between pressure measurements and node degree between nearby sensors"
Set-ParagraphText $tr 14 "e.g., landslide monitoring sensor network: correlation between pressure measurements and node degree between nearby sensors"

# 3) "LLMs:" -> "Multistore:"
Set-ParagraphText $tr 16 "Multistore:"

# 4) Drop the "Text to query (in hybrid models)" bullet (paragraph 17) entirely, and turn
#    the former "Repair" bullet (paragraph 18, same outline level) into the new bullet
#    text. Deleting paragraph 17 (not the very last paragraph of the shape) keeps the
#    removal clean; then paragraph 18 slides up to become paragraph 17.
$para17 = $tr.Paragraphs(17, 1)
$toRemove = $tr.Characters($para17.Start, $para17.Length)
$toRemove.Delete()

Set-ParagraphText $tr 17 "Provide a unified language that transparently distributes the execution plan on different engines"
